$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text before writing, so that
# numeric-looking strings (e.g. "9.00", "80.77") are not auto-coerced
# into Excel numbers -- the source data is plain text in every row.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '39.203.68'
$ws.Range("E2").Value = '  -2.06%  '

# Row 3
$ws.Range("D3").Value = '2.200.46'
$ws.Range("E3").Value = '  -5.60%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '294.56'
$ws.Range("E5").Value = '  -4.28%  '

# Row 6
$ws.Range("D6").Value = '80.77'
$ws.Range("E6").Value = '  -4.86%  '

# Row 7
$ws.Range("E7").Value = '  -4.11%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").Value = '0.466'
$ws.Range("E9").Value = '  -3.77%  '

# Row 10
$ws.Range("D10").Value = '0.0767'
$ws.Range("E10").Value = '  -5.64%  '

# Row 11
$ws.Range("D11").Value = '29.04'
$ws.Range("E11").Value = '  -3.40%  '

# Row 12
$ws.Range("D12").Value = '46.85'

# Row 13
$ws.Range("E13").Value = '  -2.75%  '

# Row 14
$ws.Range("D14").Value = '2.533.65'
$ws.Range("E14").Value = '  -5.77%  '

# Row 15
$ws.Range("D15").Value = '6.20'
$ws.Range("E15").Value = '  -3.50%  '

# Row 16
$ws.Range("D16").Value = '13.88'
$ws.Range("E16").Value = '  -5.49%  '

# Row 17
$ws.Range("D17").Value = '2.193.31'
$ws.Range("E17").Value = '  -5.78%  '

# Row 18
$ws.Range("D18").Value = '0.708'
$ws.Range("E18").Value = '  -5.98%  '

# Row 19
$ws.Range("D19").Value = '39.125.66'
$ws.Range("E19").Value = '  -2.13%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0868'
$ws.Range("E20").Value = '  -3.68%  '

# Row 21
$ws.Range("D21").Value = '5.69'
$ws.Range("E21").Value = '  -6.44%  '

# Row 22
$ws.Range("D22").Value = '64.48'
$ws.Range("E22").Value = '  -4.70%  '

# Row 23
$ws.Range("D23").Value = '10.18'
$ws.Range("E23").Value = '  -4.66%  '

# Row 24
$ws.Range("D24").Value = '226.57'
$ws.Range("E24").Value = '  -3.74%  '

# Row 26
$ws.Range("D26").Value = '2.39'
$ws.Range("E26").Value = '  -6.47%  '

# Row 27
$ws.Range("D27").Value = '1.79'
$ws.Range("E27").Value = '  -0.88%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '22.44'
$ws.Range("E28").Value = '  -4.16%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.17'
$ws.Range("E29").Value = '  -1.49%  '

# Row 30
$ws.Range("D30").Value = '9.02'
$ws.Range("E30").Value = '  -2.71%  '

# Row 31
$ws.Range("D31").Value = '149.28'
$ws.Range("E31").Value = '  -1.96%  '

# Row 32
$ws.Range("D32").Value = '31.39'
$ws.Range("E32").Value = '  -12.01%  '

# Row 33
$ws.Range("E33").Value = '  -0.22%  '

# Row 34
$ws.Range("E34").Value = '  -6.60%  '

# Row 35
$ws.Range("E35").Value = '  -4.38%  '

# Row 36
$ws.Range("D36").Value = '0.0693'
$ws.Range("E36").Value = '  -4.48%  '

# Row 37
$ws.Range("E37").Value = '  -3.60%  '

# Row 38
$ws.Range("D38").Value = '15.30'
$ws.Range("E38").Value = '  -2.93%  '

# Row 39
$ws.Range("D39").Value = '0.0957'
$ws.Range("E39").Value = '  -3.89%  '

# Row 40
$ws.Range("E40").Value = '  -4.93%  '

# Row 41
$ws.Range("D41").Value = '1.64'
$ws.Range("E41").Value = '  -3.82%  '

# Row 42
$ws.Range("D42").Value = '3.58'
$ws.Range("E42").Value = '  -6.42%  '

# Row 43
$ws.Range("D43").Value = '1.893.71'
$ws.Range("E43").Value = '  -2.57%  '

# Row 44
$ws.Range("D44").Value = '2.03'
$ws.Range("E44").Value = '  -10.24%  '

# Row 45
$ws.Range("D45").Value = '0.0258'
$ws.Range("E45").Value = '  -3.07%  '

# Row 46
$ws.Range("D46").Value = '9.00'
$ws.Range("E46").Value = '  -2.95%  '

# Row 47
$ws.Range("E47").Value = '  -8.94%  '

# Row 48
$ws.Range("E48").Value = '  -3.11%  '

# Row 49
$ws.Range("D49").Value = '71.24'
$ws.Range("E49").Value = '  +0.81%  '

# Row 50
$ws.Range("D50").Value = '2.405.17'
$ws.Range("E50").Value = '  -5.85%  '

# Row 51
$ws.Range("D51").Value = '87.17'
$ws.Range("E51").Value = '  -6.19%  '

# Remove the temporary Text number-format so the cells keep their
# original (default) style -- only the values themselves changed.
$priceVolRange.ClearFormats()
